# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Mon Jan 22 11:47:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.775.83'
$ws.Range("E2").Value = '  -2.24%  '
$ws.Range("D3").Value = '2.382.79'
$ws.Range("E3").Value = '  -3.73%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '88.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.530'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.02%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0823'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '31.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.43%  '
$ws.Range("E12").Value = '  -1.67%  '
$ws.Range("D13").Value = '2.758.05'
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.79%  '
$ws.Range("D16").Value = '2.382.86'
$ws.Range("E16").Value = '  -3.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.762'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.94%  '
$ws.Range("D18").Value = '40.753.94'
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.18%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  -4.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.47%  '
$ws.Range("E24").Value = '  -3.49%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  -6.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.81%  '
$ws.Range("E28").Value = '  -2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.89'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.69%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0732'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("E35").Value = '  -6.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.20%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.08%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.16%  '
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.71%  '
$ws.Range("E42").Value = '  -7.91%  '
$ws.Range("D43").Value = '1.955.15'
$ws.Range("E43").Value = '  -1.93%  '
$ws.Range("E44").Value = '  -5.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("D48").Value = '2.615.81'
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '73.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '93.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.99'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.44%  '
